# "Signed Off time sheets"
# The supervisor (Ankita Gangotra) has now filled in her name and signed off
# the timesheet: her name goes into the "Supervisor Name" field, and her
# initials + sign-off date go into the second signature line (mirroring the
# existing employee signature line directly above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor Name (field next to "Supervisor Name:" label, row 6)
$ws.Range("G6").Value = "Ankita Gangotra"

# Copy the formatting of the employee signature/date line (row 25) down onto
# the previously-blank supervisor signature/date line (row 27) so the new
# entries look consistent with the employee's.
$ws.Range("A25:E25").Copy() | Out-Null
$ws.Range("A27:E27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Supervisor's signature initials and sign-off date
$ws.Range("A27").Value = "A.G"
$ws.Range("D27").Value = 41800

$excel.CutCopyMode = 0
